# Update "想去人数" (interest count) values in column F across the four
# worksheets of the 广州-漫展信息 workbook, per the data refresh recorded
# in the commit "Update gh-pages to output generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 13959
$ws1.Range("F4").Value  = 13791
$ws1.Range("F12").Value = 785
$ws1.Range("F17").Value = 155
$ws1.Range("F19").Value = 559
$ws1.Range("F21").Value = 463
$ws1.Range("F23").Value = 12
$ws1.Range("F27").Value = 39
$ws1.Range("F28").Value = 9
$ws1.Range("F31").Value = 16

# --- Sheet "演出" (Performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F6").Value  = 87
$ws2.Range("F15").Value = 1653

# --- Sheet "本地生活" (Local Life) ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 33

# --- Sheet "全部类型" (All Types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 13959
$ws4.Range("F5").Value  = 13791
$ws4.Range("F13").Value = 785
$ws4.Range("F17").Value = 33
$ws4.Range("F21").Value = 155
$ws4.Range("F24").Value = 87
$ws4.Range("F26").Value = 559
$ws4.Range("F28").Value = 463
$ws4.Range("F30").Value = 12
$ws4.Range("F40").Value = 39
$ws4.Range("F41").Value = 9
$ws4.Range("F46").Value = 16
$ws4.Range("F48").Value = 1653
